# Update to 2023.2 release: rename resiliency-project categories and
# refresh the active cell selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserInputs")

# "Road_complete" -> "Road"
$ws.Range("D2").Value = "Road"

# "Subway_complete" -> "Rail"
$ws.Range("D3").Value = "Rail"

# Leave the active cell on D3, matching the saved selection in the workbook.
$ws.Range("D3").Select()
